# Update column G ("K" = strikeouts) values for rows 2-32.
# These values were regenerated (using K instead of Strike#, recalculated
# std/mean) and written back into the save_data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 7
    3  = 8
    4  = 7
    5  = 2
    6  = 5
    7  = 6
    8  = 4
    9  = 9
    10 = 4
    11 = 2
    12 = 9
    13 = 8
    14 = 11
    15 = 5
    16 = 6
    17 = 6
    18 = 4
    19 = 5
    20 = 4
    21 = 5
    22 = 5
    23 = 3
    24 = 6
    25 = 6
    26 = 5
    27 = 7
    28 = 3
    29 = 4
    30 = 2
    31 = 5
    32 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
